# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Golem_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the target diff (profit/price recalculation updates for specific leve rows).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1183
$ws.Range("I32").Value = 1183
$ws.Range("K32").Value = 1183
$ws.Range("M32").Value = -857
# Row 33
$ws.Range("H33").Value = 55.5
$ws.Range("I33").Value = 55.5
$ws.Range("K33").Value = 55.5
$ws.Range("M33").Value = 173.5
# Row 53
$ws.Range("H53").Value = 132.47058
$ws.Range("I53").Value = 127
$ws.Range("K53").Value = 127
$ws.Range("M53").Value = 510
# Row 101
$ws.Range("H101").Value = 4168
$ws.Range("I101").Value = 4168
$ws.Range("K101").Value = 12504
$ws.Range("M101").Value = -10882
# Row 105
$ws.Range("H105").Value = 21055.625
$ws.Range("J105").Value = 21055.625
$ws.Range("L105").Value = 21055.625
$ws.Range("N105").Value = -28043.625
# Row 111
$ws.Range("H111").Value = 880
$ws.Range("I111").Value = 500
$ws.Range("K111").Value = 1500
$ws.Range("M111").Value = 1567
# Row 112
$ws.Range("H112").Value = 1311.0714
$ws.Range("I112").Value = 855
$ws.Range("K112").Value = 2565
$ws.Range("M112").Value = -1457
# Row 113
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8508
# Row 129
$ws.Range("H129").Value = 18248.75
$ws.Range("I129").Value = 18248.75
$ws.Range("K129").Value = 54746.25
$ws.Range("M129").Value = -49746.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 717.3
$ws.Range("I32").Value = 523.5
$ws.Range("J32").Value = 1492.5
$ws.Range("K32").Value = 523.5
$ws.Range("L32").Value = 1492.5
$ws.Range("M32").Value = -236.5
$ws.Range("N32").Value = -2066.5
# Row 45
$ws.Range("H45").Value = 1248.5
$ws.Range("I45").Value = 1278.2
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 1278.2
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -901.2
$ws.Range("N45").Value = -1854
# Row 50
$ws.Range("H50").Value = 19842.572
$ws.Range("I50").Value = 11474
$ws.Range("K50").Value = 11474
$ws.Range("M50").Value = -10760
# Row 61
$ws.Range("H61").Value = 1200
$ws.Range("I61").Value = 1200
$ws.Range("K61").Value = 1200
$ws.Range("M61").Value = -988
# Row 74
$ws.Range("H74").Value = 1621.25
$ws.Range("I74").Value = 1495
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1495
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -621
$ws.Range("N74").Value = -3748
# Row 77
$ws.Range("H77").Value = 1621.25
$ws.Range("I77").Value = 1495
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 7475
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -3107
$ws.Range("N77").Value = -18736
# Row 133
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -45060
# Row 136
$ws.Range("H136").Value = 1200
$ws.Range("I136").Value = 1200
$ws.Range("K136").Value = 3600
$ws.Range("M136").Value = -1050

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 800
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -627
# Row 26
$ws.Range("H26").Value = 27082.572
$ws.Range("I26").Value = 27082.572
$ws.Range("K26").Value = 27082.572
$ws.Range("M26").Value = -26790.572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 149999.67
$ws.Range("J9").Value = 149999.67
$ws.Range("L9").Value = 149999.67
$ws.Range("N9").Value = -150335.67
# Row 22
$ws.Range("H22").Value = 675.4545000000001
$ws.Range("I22").Value = 675.4545000000001
$ws.Range("K22").Value = 675.4545000000001
$ws.Range("M22").Value = -325.4545000000001
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 94
$ws.Range("H94").Value = 1863.5
$ws.Range("I94").Value = 870.3333
$ws.Range("J94").Value = 2856.6667
$ws.Range("K94").Value = 870.3333
$ws.Range("L94").Value = 2856.6667
$ws.Range("M94").Value = -419.3333
$ws.Range("N94").Value = -3758.6667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1464.1666
$ws.Range("J5").Value = 1928.3334
$ws.Range("L5").Value = 5785.0002
$ws.Range("N5").Value = -6009.0002
# Row 38
$ws.Range("H38").Value = 753.63635
$ws.Range("I38").Value = 49.666668
$ws.Range("K38").Value = 149.000004
$ws.Range("M38").Value = 197.999996
# Row 59
$ws.Range("H59").Value = 1800
$ws.Range("J59").Value = 1800
$ws.Range("L59").Value = 5400
$ws.Range("N59").Value = -6480
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
# Row 109
$ws.Range("H109").Value = 932
$ws.Range("I109").Value = 942.6667
$ws.Range("K109").Value = 2828.0001
$ws.Range("M109").Value = -1788.0001
# Row 129
$ws.Range("H129").Value = 2396.3333
$ws.Range("J129").Value = 6500
$ws.Range("L129").Value = 19500
$ws.Range("N129").Value = -29500
# Row 135
$ws.Range("H135").Value = 1464.1666
$ws.Range("J135").Value = 1928.3334
$ws.Range("L135").Value = 17355.0006
$ws.Range("N135").Value = -22425.0006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 29736.285
$ws.Range("I102").Value = 34459
$ws.Range("K102").Value = 34459
$ws.Range("M102").Value = -32837
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 2757
$ws.Range("J132").Value = 2514
$ws.Range("L132").Value = 7542
$ws.Range("N132").Value = -12602

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2136.5625
$ws.Range("I22").Value = 934.25
$ws.Range("J22").Value = 3338.875
$ws.Range("K22").Value = 934.25
$ws.Range("L22").Value = 3338.875
$ws.Range("M22").Value = -639.25
$ws.Range("N22").Value = -3928.875
# Row 27
$ws.Range("H27").Value = 2136.5625
$ws.Range("I27").Value = 934.25
$ws.Range("J27").Value = 3338.875
$ws.Range("K27").Value = 934.25
$ws.Range("L27").Value = 3338.875
$ws.Range("M27").Value = -827.25
$ws.Range("N27").Value = -3552.875
# Row 46
$ws.Range("H46").Value = 3999.4443
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -812
# Row 55
$ws.Range("H55").Value = 611.1111
$ws.Range("J55").Value = 962.25
$ws.Range("L55").Value = 962.25
$ws.Range("N55").Value = -1308.25
# Row 68
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 5491
$ws.Range("I6").Value = 12252.5
$ws.Range("J6").Value = 983.3333
$ws.Range("K6").Value = 12252.5
$ws.Range("L6").Value = 983.3333
$ws.Range("M6").Value = -12137.5
$ws.Range("N6").Value = -1213.3333
# Row 113
$ws.Range("H113").Value = 5917.7
$ws.Range("I113").Value = 397.125
$ws.Range("K113").Value = 1191.375
$ws.Range("M113").Value = 978.625
# Row 121
$ws.Range("H121").Value = 78000
$ws.Range("J121").Value = 78000
$ws.Range("L121").Value = 78000
$ws.Range("N121").Value = -81494

